$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 69 entirely (the duplicate "Home/Application/Not Complete" entry
# for 45346), which shifts all subsequent rows up by one and leaves the
# 45347 entry (formerly row 70) as the new row 69.
$ws.Rows("69").Delete()

# Update the view: top-left cell, zoom, and selection as recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("F69").Select()
